$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 1

$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1

$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1

$ws.Range("M11").Value = 1

$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 1

$ws.Range("M14").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 0

$ws.Range("E21").Formula = "=(SUM(E17:E18)+SUM(E13:E15)+SUM(E9:E11)+SUM(E5:E7))+F21"
$ws.Range("F21").Formula = "=(SUM(F17:F18)+SUM(F13:F15)+SUM(F9:F11)+SUM(F5:F7))+G21"
$ws.Range("G21").Formula = "=(SUM(G17:G18)+SUM(G13:G15)+SUM(G9:G11)+SUM(G5:G7))+H21"
$ws.Range("H21").Formula = "=(SUM(H17:H18)+SUM(H13:H15)+SUM(H9:H11)+SUM(H5:H7))+I21"
$ws.Range("I21").Formula = "=(SUM(I17:I18)+SUM(I13:I15)+SUM(I9:I11)+SUM(I5:I7))+J21"
$ws.Range("J21").Formula = "=(SUM(J17:J18)+SUM(J13:J15)+SUM(J9:J11)+SUM(J5:J7))+K21"
$ws.Range("K21").Formula = "=(SUM(K17:K18)+SUM(K13:K15)+SUM(K9:K11)+SUM(K5:K7))+L21"
$ws.Range("L21").Formula = "=(SUM(L17:L18)+SUM(L13:L15)+SUM(L9:L11)+SUM(L5:L7))+M21"
$ws.Range("M21").Formula = "=(SUM(M17:M18)+SUM(M13:M15)+SUM(M9:M11)+SUM(M5:M7))+N21"
$ws.Range("N21").Formula = "=(SUM(N17:N18)+SUM(N13:N15)+SUM(N9:N11)+SUM(N5:N7))+O21"
$ws.Range("O21").Formula = "=(SUM(O17:O18)+SUM(O13:O15)+SUM(O9:O11)+SUM(O5:O7))+P21"

# View: zoom + scroll position + selection, matching the author's re-saved window state
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollColumn = 3
$win.ScrollRow = 1
[void]$ws.Range("P34").Select()
